$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date label
$ws.Range("I1").Value = "24/03/2023"

# Data for rows 2..12, columns B..J (A is label, unchanged)
$data = @{
    2  = @(14, 339, 401, 12, 0, 6, 35, 601.9, -33.37763748130919)
    3  = @(1, 46, 49, 2, 0, 3, 0, 79, -37.9746835443038)
    4  = @(2, 125, 130, 3, 1, 1, 0, 100, 30)
    5  = @(5, 281, 311, 5, 2, 5, 18, 558, -44.2652329749104)
    6  = @(30, 212, 290, 45, 4, 0, 0, 366, -20.76502732240437)
    7  = @(0, 77, 81, 4, 0, 0, 0, 115, -29.56521739130434)
    8  = @(0, 97, 124, 25, 2, 2, 0, 151, -17.88079470198676)
    9  = @(0, 131, 131, 0, 0, 2, 0, 392, -66.58163265306123)
    10 = @(0, 18, 18, 0, 0, 1, 0, 47, -61.70212765957447)
    11 = @(0, 0, 0, 0, 0, 0, 0, 1, -100)
    12 = @(0, 13, 13, 0, 0, 0, 0, 34, -61.76470588235294)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
